$wb = $excel.ActiveWorkbook

# --- Revert metadata values on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "1.1.0"
$meta.Range("B8").Value  = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Revert the descendant-of value id on the "Include from FSIII" sheet ---
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Range("C2").Value = "B"

# --- Remove the extra "Include from FSIII 2" sheet that was added ---
$extra = $wb.Worksheets.Item("Include from FSIII 2")
$extra.Delete()
